$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 115, shifting existing rows 115-144 down to 116-145.
$ws.Rows.Item(115).Insert()

# Populate the newly inserted row 115 with the new data record.
$ws.Cells.Item(115, 1).Value = 5
$ws.Cells.Item(115, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(115, 3).Value = "Maule"
$ws.Cells.Item(115, 4).Value = 45173
$ws.Cells.Item(115, 5).Value = 7
$ws.Cells.Item(115, 6).Value = 100112013
$ws.Cells.Item(115, 7).Value = "Alcachofa"
$ws.Cells.Item(115, 8).Value = "Madrigal"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 300
$ws.Cells.Item(115, 11).Value = 10000
$ws.Cells.Item(115, 12).Value = 10000
$ws.Cells.Item(115, 13).Value = 10000
$ws.Cells.Item(115, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(115, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(115, 16).Value = 250
$ws.Cells.Item(115, 17).Value = 40
$ws.Cells.Item(115, 18).Value = "Hortaliza"
